$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D and E values stay as text (matching original inlineStr cells)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.030.49'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.577.53'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.26%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.08'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.19'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.02%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.74'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.38%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.17%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.73%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.032.06'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '57.973.18'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.67'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.572.00'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.70%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.67%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '333.27'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.03'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.17'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.90%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.78'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.420'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.01'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.64%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.57'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.02%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.79'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.75%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '36.91'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.35%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.834'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.820'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.58'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '280.89'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.588'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.05%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.64'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0949'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.79'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.51%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.906.57'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.75'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.42%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.55%  '
